$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "38.214.33"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +0.61%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.096.00"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +2.88%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "229.77"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.50%  "
$ws.Range("E6").Value = "  +0.27%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "60.88"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +0.39%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("E9").Value = "  -0.24%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0845"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +3.66%  "
$ws.Range("E11").Value = "  +0.12%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "2.405.02"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +2.72%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "22.50"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +5.59%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "14.67"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +0.42%  "
$ws.Range("E15").Value = "  +6.40%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.774"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +1.22%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.092.26"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +2.69%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "38.118.04"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +0.47%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "70.35"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +0.73%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.99"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +0.56%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0835"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +1.10%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "224.45"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +0.02%  "
$ws.Range("E23").Value = "  +0.48%  "
$ws.Range("E24").Value = "  +0.49%  "
$ws.Range("E25").Value = "  +3.54%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "170.22"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +1.76%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.44"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +1.14%  "
$ws.Range("E28").Value = "  +1.47%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.02"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +0.55%  "
$ws.Range("E30").Value = "  +5.32%  "
$ws.Range("E31").Value = "  -0.11%  "
$ws.Range("E32").Value = "  +8.61%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.69"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +3.71%  "
$ws.Range("E34").Value = "  +0.50%  "
$ws.Range("E35").Value = "  -0.23%  "
$ws.Range("E36").Value = "  +0.71%  "
$ws.Range("E37").Value = "  +5.13%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.55"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +7.80%  "
$ws.Range("E39").Value = "  -0.18%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.08"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +1.98%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.550.35"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +1.50%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "100.31"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +4.17%  "
$ws.Range("E43").Value = "  +0.36%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.84"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +0.20%  "
$ws.Range("E45").Value = "  -1.24%  "
$ws.Range("E46").Value = "  +1.28%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.12"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +2.40%  "
$ws.Range("E48").Value = "  +1.67%  "
$ws.Range("E49").Value = "  +1.87%  "
$ws.Range("E50").Value = "  +0.95%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.292.89"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +2.80%  "
